# Tidsregistrering i PTE projektet - tilfoej rettet tidsplan for 20-03-2017
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- Seed styles for the new rows by copying formats from an existing,
# already-filled "day block" (row 42 has the same shape: date in col A,
# role/activity/start/end time in cols E:H) so the new cells reuse the
# existing style indices (date format, time format) instead of minting
# brand-new ones.
$ws.Range("A42").Copy($ws.Range("A49"))
$ws.Range("G42:H42").Copy($ws.Range("G49:H49"))
$ws.Range("G42:H42").Copy($ws.Range("G50:H50"))
$ws.Range("G42:H42").Copy($ws.Range("G51:H51"))
$ws.Range("G42:H42").Copy($ws.Range("G52:H52"))
$ws.Range("G42:H42").Copy($ws.Range("G53:H53"))

# Row 49 - ny blok for den 20-03-2017 (dato = 42814)
$ws.Range("A49").Value2 = 42814
$ws.Range("E49").Value2 = "Requirements Specifier"
$ws.Range("F49").Value2 = "Lavet OC11: getSigmaRef"
$ws.Range("G49").Value2 = 0.3576388888888889
$ws.Range("H49").Value2 = 0.37152777777777773

# Row 50
$ws.Range("F50").Value2 = "Lavet OC14: getBøjningsMoment"
$ws.Range("G50").Value2 = 0.37152777777777773
$ws.Range("H50").Value2 = 0.41666666666666669

# Row 51
$ws.Range("E51").Value2 = "Reviewer"
$ws.Range("F51").Value2 = "Lavet review af OC11: getSigmaRef"
$ws.Range("G51").Value2 = 0.41666666666666669
$ws.Range("H51").Value2 = 0.42708333333333331

# Row 52
$ws.Range("E52").Value2 = "Requirements Specifier"
$ws.Range("F52").Value2 = "Lavet OC9 design"
$ws.Range("G52").Value2 = 0.4375
$ws.Range("H52").Value2 = 0.54513888888888895

# Row 53
$ws.Range("F53").Value2 = "Lavet OC8 design"
$ws.Range("G53").Value2 = 0.54861111111111105
$ws.Range("H53").Value2 = 0.61805555555555558

# Row 54 - samlet tid for dagen
$ws.Range("I54").Value2 = 5.55

# Scroll/selection state, matching the author's view when they saved
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I51").Select() | Out-Null
